$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 104 originally held a broken "=TODAY==" formula instead of a file-name
# entry; removing it re-aligns the File/Review date pairing for every row
# below it (the dates in column B were already correct, just offset by one).
$ws.Rows.Item(104).Delete()

# Re-generate the File column (A) so each file name lines up with its true
# review date and the whole list reads in ascending date order.
$ws.Cells.Item(2, 1).Value = "Routine_Care/Nursing for Arterial and Central Venous Lines.pdf"
$ws.Cells.Item(3, 1).Value = "Routine_Care/VTE_Prevention/TED Stocking Sizing.pdf"
$ws.Cells.Item(4, 1).Value = "Routine_Care/Faecal  incontinence skin care.pdf"
$ws.Cells.Item(5, 1).Value = "Breathing(Respiratory)/Equipment/IPPB using an ICU Ventilator.pdf"
$ws.Cells.Item(6, 1).Value = "Drugs/heparin_critical_care_only.pdf"
$ws.Cells.Item(7, 1).Value = "Routine_Care/Invasive Flush Systems.pdf"
$ws.Cells.Item(8, 1).Value = "GI_Liver_and_Transplant/Pancreatic Irrigation.pdf"
$ws.Cells.Item(9, 1).Value = "Neurological/Thiopentone levels.pdf"
$ws.Cells.Item(10, 1).Value = "End_of_life_care/End of life care in critical care.pdf"
$ws.Cells.Item(11, 1).Value = "Infection_and_sepsis/SARI/Reporting influenza deaths.pdf"
$ws.Cells.Item(12, 1).Value = "Infection_and_sepsis/Ebola/Ebola.pdf"
$ws.Cells.Item(13, 1).Value = "GI_Liver_and_Transplant/Nasal bridle.pdf"
$ws.Cells.Item(14, 1).Value = "Infection_and_sepsis/SARI/Setup Guide for Jupiter Hoods.pdf"
$ws.Cells.Item(15, 1).Value = "Breathing(Respiratory)/Equipment/AMBU AScope.pdf"
$ws.Cells.Item(16, 1).Value = "Cardiovascular/EZ-IO Intraosseus Access Device_pub_em.pdf"
$ws.Cells.Item(17, 1).Value = "Routine_Care/Central venous catheter removal.pdf"
$ws.Cells.Item(18, 1).Value = "Routine_Care/Tracheostomy_nursing_care.pdf"
$ws.Cells.Item(19, 1).Value = "Infection_and_sepsis/SARI/MERS-CoV Guideline.pdf"
$ws.Cells.Item(20, 1).Value = "Infection_and_sepsis/SARI/Suspected Influeza A Avian influenza H5N1 and SARS.pdf"
$ws.Cells.Item(21, 1).Value = "Infection_and_sepsis/SARI/Suspected Influeza A-H7N9 Guideline.pdf"
$ws.Cells.Item(22, 1).Value = "Infection_and_sepsis/SARI/Management of Patients with severe acute respiratory infection SARI.pdf"
$ws.Cells.Item(23, 1).Value = "Infection_and_sepsis/SARI/Management of patients with SARI-additional Information.pdf"
$ws.Cells.Item(24, 1).Value = "ECLS/Extra Corporeal Carbon Dioxide Removal.pdf"
$ws.Cells.Item(25, 1).Value = "Airway/Critical care extubation checklist.pdf"
$ws.Cells.Item(26, 1).Value = "Airway/Tracheostomy_Laryngectomy/Hospital_in-patients_with_a_Tracheostomy.pdf"
$ws.Cells.Item(27, 1).Value = "Trauma and Burns/Mangement of burns.pdf"
$ws.Cells.Item(28, 1).Value = "Drugs/diazepam_diazemuls.pdf"
$ws.Cells.Item(29, 1).Value = "End_of_life_care/Reasons to report a death to PF.pdf"
$ws.Cells.Item(30, 1).Value = "Breathing(Respiratory)/salbutamol and ipratroprium MDI.pdf"
$ws.Cells.Item(31, 1).Value = "Airway/Tracheostomy_Laryngectomy/Tracheostomy change in Critical Care.pdf"
$ws.Cells.Item(32, 1).Value = "Airway/Tracheostomy_Laryngectomy/Tracheostomy suctioning cleaning guideline.pdf"
$ws.Cells.Item(33, 1).Value = "Neurological/SOP -  Femoral site care.pdf"
$ws.Cells.Item(34, 1).Value = "Drugs/sodium_bicarbonate.pdf"
$ws.Cells.Item(35, 1).Value = "Drugs/paracetamol.pdf"
$ws.Cells.Item(36, 1).Value = "Post_op_care/Anticoagulation antiplatelet agents and epidural analgesia.pdf"
$ws.Cells.Item(37, 1).Value = "Post_op_care/Epidural top-up.pdf"
$ws.Cells.Item(38, 1).Value = "Diabetes_and_Glucose/Hyperosmolar Hyperglycaemic State.pdf"
$ws.Cells.Item(39, 1).Value = "Airway/Emergency intubation checklist_em_pub.pdf"
$ws.Cells.Item(40, 1).Value = "Covid-19/SJH/SJH COVID19 ITU Intubation Action Card.pdf"
$ws.Cells.Item(41, 1).Value = "Covid-19/SJH/SJH COVID19 ED Intubation Action Card.pdf"
$ws.Cells.Item(42, 1).Value = "Covid-19/WGH/CoVid intubation checklist WGH.pdf"
$ws.Cells.Item(43, 1).Value = "Drugs/heparin for Haemofiltration.pdf"
$ws.Cells.Item(44, 1).Value = "Drugs/fentanyl.pdf"
$ws.Cells.Item(45, 1).Value = "Airway/Tracheostomy_Laryngectomy/Tracheostomy guideline.pdf"
$ws.Cells.Item(46, 1).Value = "Covid-19/WGH/WGH_CT_Transfer_May.pdf"
$ws.Cells.Item(47, 1).Value = "Cardiovascular/GJNH Acute Heart Failure Referral Form.pdf"
$ws.Cells.Item(48, 1).Value = "Organ_donation/Donation after circulatory death.pdf"
$ws.Cells.Item(49, 1).Value = "Airway/Percutaneous tracheostomy checklist.pdf"
$ws.Cells.Item(50, 1).Value = "Delirium/Managing a Potentially Violent Patient.pdf"
$ws.Cells.Item(51, 1).Value = "Delirium/Risk assessment posi mit.pdf"
$ws.Cells.Item(52, 1).Value = "Infection_and_sepsis/SOP Ultrasound Cleaning.pdf"
$ws.Cells.Item(53, 1).Value = "Breathing(Respiratory)/HFNO.pdf"
$ws.Cells.Item(54, 1).Value = "GI_Liver_and_Transplant/Treatment of constipation.pdf"
$ws.Cells.Item(55, 1).Value = "GI_Liver_and_Transplant/Abdominal pressure measurement.pdf"
$ws.Cells.Item(56, 1).Value = "Delirium/Drugs Causing Delirium and Agitiation.pdf"
$ws.Cells.Item(57, 1).Value = "Airway/McGrath Mac.pdf"
$ws.Cells.Item(58, 1).Value = "Neurological/Sub arachnoid haemorrhage management.pdf"
$ws.Cells.Item(59, 1).Value = "Airway/Tracheostomy_Laryngectomy/Tracheostomy safety box contents.pdf"
$ws.Cells.Item(60, 1).Value = "Drugs/ketamine_in_asthma.pdf"
$ws.Cells.Item(61, 1).Value = "Airway/Anticipated difficult airway tool.pdf"
$ws.Cells.Item(62, 1).Value = "End_of_life_care/Documentation following death.pdf"
$ws.Cells.Item(63, 1).Value = "Drugs/anidulafungin.pdf"
$ws.Cells.Item(64, 1).Value = "Drugs/zanamivir.pdf"
$ws.Cells.Item(65, 1).Value = "GI_Liver_and_Transplant/stress ulcer prophylaxis.pdf"
$ws.Cells.Item(66, 1).Value = "Routine_Care/bBraun Spaceplus Failure EMERGENCY ACTION CARD_em.pdf"
$ws.Cells.Item(67, 1).Value = "Drugs/phosphate.pdf"
$ws.Cells.Item(68, 1).Value = "Drugs/insulin.pdf"
$ws.Cells.Item(69, 1).Value = "Breathing(Respiratory)/Equipment/HFNO Set Up.pdf"
$ws.Cells.Item(70, 1).Value = "Breathing(Respiratory)/Inhaled Nitrous Oxide.pdf"
$ws.Cells.Item(71, 1).Value = "Breathing(Respiratory)/Equipment/APRV.pdf"
$ws.Cells.Item(72, 1).Value = "Cardiovascular/Steroids for Septic Shock.pdf"
$ws.Cells.Item(73, 1).Value = "Drugs/midazolam.pdf"
$ws.Cells.Item(74, 1).Value = "Post_op_care/Epidural Haematoma.pdf"
$ws.Cells.Item(75, 1).Value = "Breathing(Respiratory)/Equipment/T piece Y piece.pdf"
$ws.Cells.Item(76, 1).Value = "Neurological/SOP for review of Neurosurgical patients in ITU by neurosurgical team.pdf"
$ws.Cells.Item(77, 1).Value = "Drugs/morphine.pdf"
$ws.Cells.Item(78, 1).Value = "Policies_and_admin/General Critical Care Interaction with HEPMA_pub.pdf"
$ws.Cells.Item(79, 1).Value = "Drugs/atracurium.pdf"
$ws.Cells.Item(80, 1).Value = "Infection_and_sepsis/BAL and MiniBAL standardised procedure.pdf"
$ws.Cells.Item(81, 1).Value = "Drugs/midazolam and thiopental levels.pdf"
$ws.Cells.Item(82, 1).Value = "Drugs/all IV drug infusion information.pdf"
$ws.Cells.Item(83, 1).Value = "Drugs/epoprostenol.pdf"
$ws.Cells.Item(84, 1).Value = "Breathing(Respiratory)/Equipment/HFNO through ventilator.pdf"
$ws.Cells.Item(85, 1).Value = "Routine_Care/VTE_Prevention/Dalteparin_thromboprophylaxis.pdf"
$ws.Cells.Item(86, 1).Value = "Post_op_care/Adult Scoliosis Spinal Surgery Post-Op Care.pdf"
$ws.Cells.Item(87, 1).Value = "Post_op_care/Post op care pharyngo-laryngo-oesphagectomy PLOG.pdf"
$ws.Cells.Item(88, 1).Value = "Drugs/amiodarone.pdf"
$ws.Cells.Item(89, 1).Value = "Drugs/phenylephrine.pdf"
$ws.Cells.Item(90, 1).Value = "Infection_and_sepsis/SARI/Influenza A and B Virology Sampling and Oseltamivir Dose.pdf"
$ws.Cells.Item(91, 1).Value = "GI_Liver_and_Transplant/Nasogastric feeding protocol.pdf"
$ws.Cells.Item(92, 1).Value = "Drugs/potassium.pdf"
$ws.Cells.Item(93, 1).Value = "Cardiovascular/Central Venous Catheters - Guideline for Management of Misplacement.pdf"
$ws.Cells.Item(94, 1).Value = "Renal_and_Urology/Antibiotic doses in CVVHD.pdf"
$ws.Cells.Item(95, 1).Value = "Neurological/Intrathecal policy RIE.pdf"
$ws.Cells.Item(96, 1).Value = "Diabetes_and_Glucose/Intravenous Insulin Therapy (not for DKA or HHS).pdf"
$ws.Cells.Item(97, 1).Value = "Infection_and_sepsis/Antibiotic doses in CVVHD.pdf"
$ws.Cells.Item(98, 1).Value = "Drugs/Antibiotic doses in CVVHD.pdf"
$ws.Cells.Item(99, 1).Value = "Drugs/nimodipine.pdf"
$ws.Cells.Item(100, 1).Value = "ECLS/RIE ECLS Anti Xa Protocol.pdf"
$ws.Cells.Item(101, 1).Value = "GI_Liver_and_Transplant/Nasojejunal feeding protocol.pdf"
$ws.Cells.Item(102, 1).Value = "Drugs/calcium.pdf"
$ws.Cells.Item(103, 1).Value = "GI_Liver_and_Transplant/Jejunostomy feeding protocol.pdf"
$ws.Cells.Item(104, 1).Value = "Covid-19/COVID 19 ICM guidance basic goals_June_2022.pdf"
$ws.Cells.Item(105, 1).Value = "End_of_life_care/CMO & NRS Guidance for Doctors completing MCCD - Sept 22.pdf"
$ws.Cells.Item(106, 1).Value = "Ethics_and_Law/Care at the End of Life (FICM).pdf"
$ws.Cells.Item(107, 1).Value = "Drugs/vasopressin organ donation.pdf"
$ws.Cells.Item(108, 1).Value = "Neurological/Management of traumatic brain injury.pdf"
$ws.Cells.Item(109, 1).Value = "Infection_and_sepsis/Winter Infections Stepdown Guidance.pdf"
$ws.Cells.Item(110, 1).Value = "Neurological/Critical Care MRI Procedure_pub.pdf"
$ws.Cells.Item(111, 1).Value = "Drugs/vasopressin_sepsis.pdf"
$ws.Cells.Item(112, 1).Value = "Drugs/nicardipine.pdf"
$ws.Cells.Item(113, 1).Value = "Organ_donation/Organ Retrieval SOP.pdf"
$ws.Cells.Item(114, 1).Value = "Ethics_and_Law/DNACPR policy for Scotland.pdf"
$ws.Cells.Item(115, 1).Value = "Airway/Cook Staged Extubation Set.pdf"
$ws.Cells.Item(116, 1).Value = "Post_op_care/Epidural hypotension.pdf"
$ws.Cells.Item(117, 1).Value = "Transfer/ACCP Transfers.pdf"
$ws.Cells.Item(118, 1).Value = "Drugs/ketamine for status epilepticus.pdf"
$ws.Cells.Item(119, 1).Value = "Drugs/valproate.pdf"
$ws.Cells.Item(120, 1).Value = "Drugs/thiopentone.pdf"
$ws.Cells.Item(121, 1).Value = "Covid-19/videos/Donning and Doffing Video.pdf"
$ws.Cells.Item(122, 1).Value = "Breathing(Respiratory)/Equipment/Ventilators Circuits Filters and Closed Suction - Set up and Maintenance.pdf"
$ws.Cells.Item(123, 1).Value = "Drugs/piperacillin_tazobactam extended_infusion.pdf"
$ws.Cells.Item(124, 1).Value = "Breathing(Respiratory)/Equipment/Bipap V60.pdf"
$ws.Cells.Item(125, 1).Value = "Infection_and_sepsis/Infection indications for IVIG.pdf"
$ws.Cells.Item(126, 1).Value = "Breathing(Respiratory)/CPAP.pdf"
$ws.Cells.Item(127, 1).Value = "Breathing(Respiratory)/ARDS Strategy.pdf"
$ws.Cells.Item(128, 1).Value = "Procedures/Securing CVCs.pdf"
$ws.Cells.Item(129, 1).Value = "Covid-19/Covid 19 Death Certification Guideline.pdf"
$ws.Cells.Item(130, 1).Value = "Transfer/Transfer Outdoors to Garden Guideline.pdf"
$ws.Cells.Item(131, 1).Value = "Neurological/Treatment of status epilepticus.pdf"
$ws.Cells.Item(132, 1).Value = "Routine_Care/Video Communication.pdf"
$ws.Cells.Item(133, 1).Value = "Drugs/hydralazine.pdf"
$ws.Cells.Item(134, 1).Value = "Cardiovascular/Cardiogenic Shock.pdf"
$ws.Cells.Item(135, 1).Value = "Drugs/isoprenaline.pdf"
$ws.Cells.Item(136, 1).Value = "Post_op_care/Major OMFS Free Flap.pdf"
$ws.Cells.Item(137, 1).Value = "Drugs/alfentanil.pdf"
$ws.Cells.Item(138, 1).Value = "Drugs/magnesium.pdf"
$ws.Cells.Item(139, 1).Value = "Haematology_CAR-T/Haem_ICU_transfer.pdf"
$ws.Cells.Item(140, 1).Value = "Drugs/aminophylline.pdf"
$ws.Cells.Item(141, 1).Value = "Cardiovascular/Management of hypertension within Critical Care.pdf"
$ws.Cells.Item(142, 1).Value = "Drugs/rocuronium.pdf"
$ws.Cells.Item(143, 1).Value = "Drugs/phenytoin.pdf"
$ws.Cells.Item(144, 1).Value = "Haematology_CAR-T/ICANS.pdf"
$ws.Cells.Item(145, 1).Value = "Drugs/pancuronium.pdf"
$ws.Cells.Item(146, 1).Value = "Haematology_CAR-T/CRS.pdf"
$ws.Cells.Item(147, 1).Value = "Policies_and_admin/General Critical Care SOP_pub.pdf"
$ws.Cells.Item(148, 1).Value = "Drugs/Milrinone.pdf"
$ws.Cells.Item(149, 1).Value = "Drugs/clonidine.pdf"
$ws.Cells.Item(150, 1).Value = "Drugs/noradrenaline.pdf"
$ws.Cells.Item(151, 1).Value = "GI_Liver_and_Transplant/Fulminant Liver Failure.pdf"
$ws.Cells.Item(152, 1).Value = "Drugs/glyceryl_trinitrate.pdf"
$ws.Cells.Item(153, 1).Value = "Breathing(Respiratory)/Equipment/Passy Muir Valve.pdf"
$ws.Cells.Item(154, 1).Value = "Drugs/dexmedetomidine.pdf"
$ws.Cells.Item(155, 1).Value = "GI_Liver_and_Transplant/Confirmation of Nasogastric Tube Position.pdf"
$ws.Cells.Item(156, 1).Value = "Cardiovascular/Intra Aortic Balloon Pump Bedside Checks_pub.pdf"
$ws.Cells.Item(157, 1).Value = "Cardiovascular/Intra Aortic Balloon Pump Guideline_pub.pdf"
$ws.Cells.Item(158, 1).Value = "Drugs/adrenaline.pdf"
$ws.Cells.Item(159, 1).Value = "Drugs/dobutamine.pdf"
$ws.Cells.Item(160, 1).Value = "Drugs/vancomycin.pdf"
$ws.Cells.Item(161, 1).Value = "Drugs/neostigmine.pdf"
$ws.Cells.Item(162, 1).Value = "Infection_and_sepsis/Initial investigation and management in unidentified Infections.pdf"
$ws.Cells.Item(163, 1).Value = "Drugs/labetalol.pdf"
$ws.Cells.Item(164, 1).Value = "Drugs/alteplase for massive PE.pdf"
$ws.Cells.Item(165, 1).Value = "Cardiovascular/Management of Acute Type B Aortic Dissection Guideline.pdf"
$ws.Cells.Item(166, 1).Value = "Drugs/salbutamol.pdf"
$ws.Cells.Item(167, 1).Value = "Drugs/phenobarbitone.pdf"
$ws.Cells.Item(168, 1).Value = "Procedures/Arterial Line insertion for ACCPs.pdf"
$ws.Cells.Item(169, 1).Value = "Routine_Care/ICU Eye Care Guideline.pdf"
$ws.Cells.Item(170, 1).Value = "Neurological/Ventriculitis Guideline.pdf"
$ws.Cells.Item(171, 1).Value = "Infection_and_sepsis/Ventriculitis.pdf"
$ws.Cells.Item(172, 1).Value = "Cardiovascular/Cardiac Output Monitoring _pub.pdf"
$ws.Cells.Item(173, 1).Value = "Cardiovascular/Pulmonary_Embolism_and_DVT/Catheter directed thrombolysis of iliofemoral DVT alteplase_pub.pdf"
$ws.Cells.Item(174, 1).Value = "Airway/Tracheostomy_Laryngectomy/Decannulation Guidline.pdf"
$ws.Cells.Item(175, 1).Value = "Procedures/Inadvertent Catheter Placement Guideline.pdf"
$ws.Cells.Item(176, 1).Value = "GI_Liver_and_Transplant/Upper GI bleeding  (Endoscopy) guideline for critical care.pdf"
$ws.Cells.Item(177, 1).Value = "Breathing(Respiratory)/Proning Guideline.pdf"
$ws.Cells.Item(178, 1).Value = "Procedures/ACCP CVC placement following completion of initial competencies.pdf"
$ws.Cells.Item(179, 1).Value = "Procedures/ACCPs acquiring initial CVC competencies.pdf"
$ws.Cells.Item(180, 1).Value = "Post_op_care/Prevention and treatment of paraplegia after major aortic procedures.pdf"
$ws.Cells.Item(181, 1).Value = "Breathing(Respiratory)/Equipment/NIV Set up in Critical Care.pdf"
$ws.Cells.Item(182, 1).Value = "Transfer/Transfer Guidelines.pdf"
